$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the existing row 163 (before the old row 164),
# shifting the old rows 164:173 down to 166:175.
$ws.Rows("164:165").Insert()

# New row 164: weekly price entry for "Flame Seedless" grapes ($/bandeja 10 kilos)
$ws.Range("A164").Value = 2
$ws.Range("B164").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 44931
$ws.Range("E164").Value = 4
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100109
$ws.Range("H164").Value = "Uva"
$ws.Range("I164").Value = 100109001
$ws.Range("J164").Value = "Uva"
$ws.Range("K164").Value = "Flame Seedless"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 1000
$ws.Range("N164").Value = 7000
$ws.Range("O164").Value = 8000
$ws.Range("P164").Value = 7500
$ws.Range("Q164").Value = "$/bandeja 10 kilos"
$ws.Range("R164").Value = "Provincia de Limarí"
$ws.Range("S164").Value = 750
$ws.Range("T164").Value = 10

# New row 165: weekly price entry for "Superior Seedless" grapes ($/bandeja 10 kilos)
$ws.Range("A165").Value = 2
$ws.Range("B165").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44931
$ws.Range("E165").Value = 4
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100109
$ws.Range("H165").Value = "Uva"
$ws.Range("I165").Value = 100109001
$ws.Range("J165").Value = "Uva"
$ws.Range("K165").Value = "Superior Seedless"
$ws.Range("L165").Value = "Primera"
$ws.Range("M165").Value = 700
$ws.Range("N165").Value = 9000
$ws.Range("O165").Value = 10000
$ws.Range("P165").Value = 9500
$ws.Range("Q165").Value = "$/bandeja 10 kilos"
$ws.Range("R165").Value = "Provincia de Limarí"
$ws.Range("S165").Value = 950
$ws.Range("T165").Value = 10
